$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 195
$ws.Range("I2").Value = 493
$ws.Range("J2").Value = 2244
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 610
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = 388
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 30
$ws.Range("S2").Value = 219
$ws.Range("T2").Value = 386
$ws.Range("U2").Value = 33
$ws.Range("V2").Value = 3367
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3433
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 18
